# Insert a new data row at row 131 (pushing existing rows 131-144 down to 132-145)
# and populate it with a new "Ajo" (garlic) price record for Feria Lagunitas de
# Puerto Montt, matching the rest of the rows in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Insert()

$ws.Range("A131").Value = 4
$ws.Range("B131").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value = "Los Lagos"
$ws.Range("D131").Value = 44449
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = 100112003
$ws.Range("G131").Value = "Ajo"
$ws.Range("H131").Value = "Chino"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 220
$ws.Range("K131").Value = 18000
$ws.Range("L131").Value = 18000
$ws.Range("M131").Value = 18000
$ws.Range("N131").Value = "$/caja 10 kilos"
$ws.Range("O131").Value = "China"
$ws.Range("P131").Value = 1800
$ws.Range("Q131").Value = 10
$ws.Range("R131").Value = "Hortaliza"
